$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.670.09"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "1.633.21"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'213.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "'0.0624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'19.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").Value = "1.860.98"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "1.636.27"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "26.659.24"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "'63.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "'210.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").Value = "'6.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").Value = "'147.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "'6.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.97%  "
$ws.Range("D29").Value = "'15.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  +4.49%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").Value = "1.168.93"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").Value = "'0.812"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'0.504"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'0.795"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "'5.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "1.770.30"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").Value = "'92.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "'54.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'7.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.32%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  -0.02%  "
